$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows of data (A5:B6)
$ws.Range("A5").Value = "dipak"
$ws.Range("B5").Value = 35133135

$ws.Range("A6").Value = "jpaslkrrj"
$ws.Range("B6").Value = 88.22654

# Auto-fit column B to the new numeric content (produces <cols> width entry)
$ws.Columns("B:B").AutoFit()

# Match the final selection left behind by the edit
$ws.Range("B5").Select()
